$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 51.59157666666666
$ws.Range("H2").Value = 154.77473
$ws.Range("I2").Value = 0.2641250550177587
$ws.Range("J2").Value = 0.2641250550177588
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.19524
$ws.Range("N2").Value = 33.58572
$ws.Range("O2").Value = 0.07847249539938134
$ws.Range("P2").Value = 0.07847249539938135
$ws.Range("Q2").Value = 577.5800827617333
$ws.Range("R2").Value = 5198.220744855599
$ws.Range("S2").Value = 0.02072655216474242
$ws.Range("T2").Value = 0.02072655216474242

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 51.59157666666666
$ws.Range("H3").Value = 154.77473
$ws.Range("I3").Value = 0.2641250550177587
$ws.Range("J3").Value = 0.2641250550177588
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.07813833333333
$ws.Range("N3").Value = 93.234415
$ws.Range("O3").Value = 0.2178407133195748
$ws.Range("P3").Value = 0.2178407133195749
$ws.Range("Q3").Value = 1603.370156481439
$ws.Range("R3").Value = 14430.33140833295
$ws.Range("S3").Value = 0.0575371903906405
$ws.Range("T3").Value = 0.05753719039064053

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 51.59157666666666
$ws.Range("H4").Value = 154.77473
$ws.Range("I4").Value = 0.2641250550177587
$ws.Range("J4").Value = 0.2641250550177588
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 91.56894199999999
$ws.Range("N4").Value = 274.706826
$ws.Range("O4").Value = 0.6418480872068143
$ws.Range("P4").Value = 0.6418480872068144
$ws.Range("Q4").Value = 4724.186091478552
$ws.Range("R4").Value = 42517.67482330697
$ws.Range("S4").Value = 0.169528161346543
$ws.Range("T4").Value = 0.1695281613465431

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 51.59157666666666
$ws.Range("H5").Value = 154.77473
$ws.Range("I5").Value = 0.2641250550177587
$ws.Range("J5").Value = 0.2641250550177588
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.822188333333335
$ws.Range("N5").Value = 26.466565
$ws.Range("O5").Value = 0.06183870407422939
$ws.Range("P5").Value = 0.06183870407422939
$ws.Range("Q5").Value = 455.1506057669389
$ws.Range("R5").Value = 4096.35545190245
$ws.Range("S5").Value = 0.01633315111583274
$ws.Range("T5").Value = 0.01633315111583274

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09891538535728452
$ws.Range("J6").Value = 0.09891538535728453
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.19524
$ws.Range("N6").Value = 33.58572
$ws.Range("O6").Value = 0.07847249539938134
$ws.Range("P6").Value = 0.07847249539938135
$ws.Range("Q6").Value = 216.3049486434667
$ws.Range("R6").Value = 1946.744537791201
$ws.Range("S6").Value = 0.007762137122377541
$ws.Range("T6").Value = 0.007762137122377544

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09891538535728452
$ws.Range("J7").Value = 0.09891538535728453
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 31.07813833333333
$ws.Range("N7").Value = 93.234415
$ws.Range("O7").Value = 0.2178407133195748
$ws.Range("P7").Value = 0.2178407133195749
$ws.Range("Q7").Value = 600.4654760528779
$ws.Range("R7").Value = 5404.189284475901
$ws.Range("S7").Value = 0.02154779810451148
$ws.Range("T7").Value = 0.02154779810451149

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09891538535728452
$ws.Range("J8").Value = 0.09891538535728453
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 91.56894199999999
$ws.Range("N8").Value = 274.706826
$ws.Range("O8").Value = 0.6418480872068143
$ws.Range("P8").Value = 0.6418480872068144
$ws.Range("Q8").Value = 1769.217568953107
$ws.Range("R8").Value = 15922.95812057796
$ws.Range("S8").Value = 0.063488650886898
$ws.Range("T8").Value = 0.06348865088689802

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09891538535728452
$ws.Range("J9").Value = 0.09891538535728453
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.822188333333335
$ws.Range("N9").Value = 26.466565
$ws.Range("O9").Value = 0.06183870407422939
$ws.Range("P9").Value = 0.06183870407422939
$ws.Range("Q9").Value = 170.4548535238779
$ws.Range("R9").Value = 1534.0936817149
$ws.Range("S9").Value = 0.00611679924349748
$ws.Range("T9").Value = 0.006116799243497482

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 112.3724673333333
$ws.Range("H10").Value = 337.117402
$ws.Range("I10").Value = 0.5752951554216499
$ws.Range("J10").Value = 0.57529515542165
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 11.19524
$ws.Range("N10").Value = 33.58572
$ws.Range("O10").Value = 0.07847249539938134
$ws.Range("P10").Value = 0.07847249539938135
$ws.Range("Q10").Value = 1258.036741188826
$ws.Range("R10").Value = 11322.33067069944
$ws.Range("S10").Value = 0.04514484643711179
$ws.Range("T10").Value = 0.04514484643711181

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 112.3724673333333
$ws.Range("H11").Value = 337.117402
$ws.Range("I11").Value = 0.5752951554216499
$ws.Range("J11").Value = 0.57529515542165
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 31.07813833333333
$ws.Range("N11").Value = 93.234415
$ws.Range("O11").Value = 0.2178407133195748
$ws.Range("P11").Value = 0.2178407133195749
$ws.Range("Q11").Value = 3492.327084643314
$ws.Range("R11").Value = 31430.94376178983
$ws.Range("S11").Value = 0.1253227070263479
$ws.Range("T11").Value = 0.1253227070263479

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 112.3724673333333
$ws.Range("H12").Value = 337.117402
$ws.Range("I12").Value = 0.5752951554216499
$ws.Range("J12").Value = 0.57529515542165
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 91.56894199999999
$ws.Range("N12").Value = 274.706826
$ws.Range("O12").Value = 0.6418480872068143
$ws.Range("P12").Value = 0.6418480872068144
$ws.Range("Q12").Value = 10289.82794364289
$ws.Range("R12").Value = 92608.45149278604
$ws.Range("S12").Value = 0.3692520950867329
$ws.Range("T12").Value = 0.3692520950867331

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 112.3724673333333
$ws.Range("H13").Value = 337.117402
$ws.Range("I13").Value = 0.5752951554216499
$ws.Range("J13").Value = 0.57529515542165
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 8.822188333333335
$ws.Range("N13").Value = 26.466565
$ws.Range("O13").Value = 0.06183870407422939
$ws.Range("P13").Value = 0.06183870407422939
$ws.Range("Q13").Value = 991.3710702960145
$ws.Range("R13").Value = 8922.33963266413
$ws.Range("S13").Value = 0.03557550687145721
$ws.Range("T13").Value = 0.03557550687145722

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 12.044915
$ws.Range("H14").Value = 36.134745
$ws.Range("I14").Value = 0.06166440420330686
$ws.Range("J14").Value = 0.06166440420330688
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 11.19524
$ws.Range("N14").Value = 33.58572
$ws.Range("O14").Value = 0.07847249539938134
$ws.Range("P14").Value = 0.07847249539938135
$ws.Range("Q14").Value = 134.8457142046
$ws.Range("R14").Value = 1213.6114278414
$ws.Range("S14").Value = 0.004838959675149589
$ws.Range("T14").Value = 0.004838959675149591

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 12.044915
$ws.Range("H15").Value = 36.134745
$ws.Range("I15").Value = 0.06166440420330686
$ws.Range("J15").Value = 0.06166440420330688
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 31.07813833333333
$ws.Range("N15").Value = 93.234415
$ws.Range("O15").Value = 0.2178407133195748
$ws.Range("P15").Value = 0.2178407133195749
$ws.Range("Q15").Value = 374.3335345832417
$ws.Range("R15").Value = 3369.001811249175
$ws.Range("S15").Value = 0.01343301779807495
$ws.Range("T15").Value = 0.01343301779807496

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 12.044915
$ws.Range("H16").Value = 36.134745
$ws.Range("I16").Value = 0.06166440420330686
$ws.Range("J16").Value = 0.06166440420330688
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 91.56894199999999
$ws.Range("N16").Value = 274.706826
$ws.Range("O16").Value = 0.6418480872068143
$ws.Range("P16").Value = 0.6418480872068144
$ws.Range("Q16").Value = 1102.94012302993
$ws.Range("R16").Value = 9926.46110726937
$ws.Range("S16").Value = 0.03957917988664035
$ws.Range("T16").Value = 0.03957917988664036

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 12.044915
$ws.Range("H17").Value = 36.134745
$ws.Range("I17").Value = 0.06166440420330686
$ws.Range("J17").Value = 0.06166440420330688
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 8.822188333333335
$ws.Range("N17").Value = 26.466565
$ws.Range("O17").Value = 0.06183870407422939
$ws.Range("P17").Value = 0.06183870407422939
$ws.Range("Q17").Value = 106.2625085889917
$ws.Range("R17").Value = 956.3625773009252
$ws.Range("S17").Value = 0.00381324684344196
$ws.Range("T17").Value = 0.003813246843441961

